$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sorted "Periodo Mora" (col E), "Valor Mora" (col F) and "Salario Basico" (col G)
# for rows 16-24 (ascending period order), per updated EC database.
$periods = @("1607","1608","1610","1611","1701","1702","1705","1802","1807")
$valorMora = @(27578,27578,27578,27578,27578,27578,27578,29509,31249)
$salarioBasico = @(877803,877803,877803,877803,877803,877803,877803,877803,877803)

for ($i = 0; $i -lt 9; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valorMora[$i]
    $ws.Cells.Item($row, 7).Value = $salarioBasico[$i]
}
